$wb = $excel.ActiveWorkbook

# Row 9 (item id 5487) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 144.11765
$ws.Range("I9").Value = 112.75
$ws.Range("J9").Value = 153.76923
$ws.Range("K9").Value = 112.75
$ws.Range("L9").Value = 153.76923
$ws.Range("M9").Value = 56.25
$ws.Range("N9").Value = -491.76923

# Row 100 (item id 19906) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 20836076
$ws.Range("I100").Value = 33334824
$ws.Range("J100").Value = 4831.1665
$ws.Range("K100").Value = 33334824
$ws.Range("L100").Value = 4831.1665
$ws.Range("M100").Value = -33334283
$ws.Range("N100").Value = -5913.1665

# Row 112 (item id 27960) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1623.2916
$ws.Range("I112").Value = 477.5
$ws.Range("J112").Value = 1852.45
$ws.Range("K112").Value = 1432.5
$ws.Range("L112").Value = 5557.35
$ws.Range("M112").Value = -324.5
$ws.Range("N112").Value = -7773.35

# Row 137 (item id 44013) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1493.6
$ws.Range("I137").Value = 1181.4286
$ws.Range("K137").Value = 3544.2858
$ws.Range("M137").Value = -994.2857999999997

# Row 138 (item id 44169) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2720.0303
$ws.Range("I138").Value = 1565
$ws.Range("J138").Value = 4287.5713
$ws.Range("K138").Value = 4695
$ws.Range("L138").Value = 12862.7139
$ws.Range("M138").Value = 445
$ws.Range("N138").Value = -23142.7139

# Row 63 (item id 12528) on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4342.5
$ws.Range("I63").Value = 1935
$ws.Range("J63").Value = 6750
$ws.Range("K63").Value = 1935
$ws.Range("L63").Value = 6750
$ws.Range("M63").Value = -1249
$ws.Range("N63").Value = -8122

# Row 66 (item id 12528) on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4342.5
$ws.Range("I66").Value = 1935
$ws.Range("J66").Value = 6750
$ws.Range("K66").Value = 9675
$ws.Range("L66").Value = 33750
$ws.Range("M66").Value = -6243
$ws.Range("N66").Value = -40614

# Row 132 (item id 43997) on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 29442928
$ws.Range("I132").Value = 40001430
$ws.Range("J132").Value = 113747.336
$ws.Range("K132").Value = 120004290
$ws.Range("L132").Value = 341242.008
$ws.Range("M132").Value = -120001760
$ws.Range("N132").Value = -346302.008

# Row 130 (item id 34682) on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 54927.8
$ws.Range("J130").Value = 54927.8
$ws.Range("L130").Value = 54927.8
$ws.Range("N130").Value = -64967.8

# Row 58 (item id 44021) on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1744.6818
$ws.Range("I58").Value = 1545.2354
$ws.Range("J58").Value = 2422.8
$ws.Range("K58").Value = 1545.2354
$ws.Range("L58").Value = 2422.8
$ws.Range("M58").Value = -1342.2354
$ws.Range("N58").Value = -2828.8

# Row 134 (item id 44020) on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9718
$ws.Range("I134").Value = 2828.5715
$ws.Range("J134").Value = 20434.889
$ws.Range("K134").Value = 8485.7145
$ws.Range("L134").Value = 61304.667
$ws.Range("M134").Value = -5950.7145
$ws.Range("N134").Value = -66374.667

# Row 136 (item id 44021) on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1744.6818
$ws.Range("I136").Value = 1545.2354
$ws.Range("J136").Value = 2422.8
$ws.Range("K136").Value = 4635.706200000001
$ws.Range("L136").Value = 7268.400000000001
$ws.Range("M136").Value = -2085.706200000001
$ws.Range("N136").Value = -12368.4

# Row 48 (item id 4724) on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2559.6
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 2559.6
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 7678.799999999999
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -8178.799999999999

# Row 109 (item id 27854) on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1164.5555
$ws.Range("I109").Value = 996.8333
$ws.Range("J109").Value = 1500
$ws.Range("K109").Value = 2990.4999
$ws.Range("L109").Value = 4500
$ws.Range("M109").Value = -1950.4999
$ws.Range("N109").Value = -6580

# Row 113 (item id 27843) on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 946.2967
$ws.Range("I113").Value = 767
$ws.Range("K113").Value = 2301
$ws.Range("M113").Value = -131

# Row 118 (item id 27872) on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 2657.3914
$ws.Range("I118").Value = 2020
$ws.Range("J118").Value = 2771.795
$ws.Range("K118").Value = 6060
$ws.Range("L118").Value = 8315.385
$ws.Range("M118").Value = -4817
$ws.Range("N118").Value = -10801.385

# Row 131 (item id 36060) on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6850331.5
$ws.Range("I131").Value = 2223.3333
$ws.Range("J131").Value = 7463595
$ws.Range("K131").Value = 6669.999899999999
$ws.Range("L131").Value = 22390785
$ws.Range("M131").Value = -1629.999899999999
$ws.Range("N131").Value = -22400865

# Row 70 (item id 14146) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20095662
$ws.Range("I70").Value = 40183624
$ws.Range("J70").Value = 7700.643
$ws.Range("K70").Value = 40183624
$ws.Range("L70").Value = 7700.643
$ws.Range("M70").Value = -40183354
$ws.Range("N70").Value = -8240.643

# Row 73 (item id 14146) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 20095662
$ws.Range("I73").Value = 40183624
$ws.Range("J73").Value = 7700.643
$ws.Range("K73").Value = 40183624
$ws.Range("L73").Value = 7700.643
$ws.Range("M73").Value = -40182688
$ws.Range("N73").Value = -9572.643

# Row 132 (item id 44008) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 28486.316
$ws.Range("I132").Value = 1816.138
$ws.Range("J132").Value = 114423.555
$ws.Range("K132").Value = 5448.414
$ws.Range("L132").Value = 343270.665
$ws.Range("M132").Value = -2918.414
$ws.Range("N132").Value = -348330.665

# Row 46 (item id 5282) on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8547.929
$ws.Range("I46").Value = 1813.5
$ws.Range("J46").Value = 13598.75
$ws.Range("K46").Value = 1813.5
$ws.Range("L46").Value = 13598.75
$ws.Range("M46").Value = -1625.5
$ws.Range("N46").Value = -13974.75

# Row 136 (item id 44060) on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 66678452
$ws.Range("I136").Value = 26951
$ws.Range("J136").Value = 90915360
$ws.Range("K136").Value = 80853
$ws.Range("L136").Value = 272746080
$ws.Range("M136").Value = -78303
$ws.Range("N136").Value = -272751180

# Row 62 (item id 12589) on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5422.6665
$ws.Range("I62").Value = 5417.3335
$ws.Range("J62").Value = 5433.3335
$ws.Range("K62").Value = 5417.3335
$ws.Range("L62").Value = 5433.3335
$ws.Range("M62").Value = -4793.3335
$ws.Range("N62").Value = -6681.3335

# Row 65 (item id 12589) on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5422.6665
$ws.Range("I65").Value = 5417.3335
$ws.Range("J65").Value = 5433.3335
$ws.Range("K65").Value = 27086.6675
$ws.Range("L65").Value = 27166.6675
$ws.Range("M65").Value = -23966.6675
$ws.Range("N65").Value = -33406.6675

# Row 136 (item id 44031) on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 36651.105
$ws.Range("I136").Value = 46297.09
$ws.Range("J136").Value = 1282.5
$ws.Range("K136").Value = 138891.27
$ws.Range("L136").Value = 3847.5
$ws.Range("M136").Value = -136341.27
$ws.Range("N136").Value = -8947.5
